$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ("H 72") was removed from the source data entirely; all rows
# below it shift up by one (row 3 -> row 2, row 4 -> row 3, ..., row 63 -> row 62).
$ws.Rows.Item(2).Delete()
